$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 44 ("Perfection" / 45119) -
# this pushes the existing rows 44-46 down to 45-47.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly record.
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 45147
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112022
$ws.Range("G44").Value = "Arveja Verde"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 50
$ws.Range("K44").Value = 26000
$ws.Range("L44").Value = 26000
$ws.Range("M44").Value = 26000
$ws.Range("N44").Value = "`$/malla 25 kilos"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 1040
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
